$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/percentage updates (Coin names, Links, Volume%) - safe as direct .Value assignment
$updates = @{
    'E2' = '  -1.71%  '
    'E3' = '  -1.76%  '
    'E4' = '  -0.07%  '
    'E5' = '  -4.33%  '
    'E6' = '  +0.19%  '
    'E7' = '  -0.13%  '
    'E8' = '  -3.68%  '
    'E9' = '  -6.97%  '
    'E10' = '  -0.22%  '
    'E11' = '  -0.36%  '
    'E12' = '  -3.74%  '
    'B13' = 'WrappedEther'
    'C13' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'E13' = '  +1.42%  '
    'B14' = 'Polkadot'
    'C14' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'E14' = '  +0.19%  '
    'E15' = '  -2.31%  '
    'E16' = '  +1.43%  '
    'E17' = '  -1.71%  '
    'E18' = '  -2.51%  '
    'E19' = '  -2.84%  '
    'B21' = 'WrappedliquidstakedEther2.0'
    'C21' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'E21' = '  -1.06%  '
    'B22' = 'Dai'
    'C22' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'E22' = '  -0.06%  '
    'E23' = '  -1.77%  '
    'E24' = '  -0.09%  '
    'E25' = '  -5.35%  '
    'E26' = '  -1.19%  '
    'E27' = '  -3.54%  '
    'E28' = '  -1.81%  '
    'E29' = '  -6.38%  '
    'E30' = '  +4.77%  '
    'E31' = '  -0.31%  '
    'E32' = '  +2.68%  '
    'E33' = '  -1.76%  '
    'E35' = '  -4.10%  '
    'E36' = '  +0.49%  '
    'E37' = '  +0.49%  '
    'E38' = '  -4.08%  '
    'E39' = '  -2.30%  '
    'E40' = '  -1.78%  '
    'E41' = '  +11.43%  '
    'E42' = '  -1.76%  '
    'E43' = '  -2.71%  '
    'E44' = '  -1.36%  '
    'E45' = '  -1.29%  '
    'E46' = '  -0.12%  '
    'E47' = '  -2.53%  '
    'E48' = '  -0.95%  '
    'E49' = '  -0.38%  '
    'E50' = '  -2.79%  '
    'E51' = '  -3.44%  '
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# Price column (D) updates - these are numeric-looking text values (e.g. "1.0000", "29.809.37")
# that must stay as literal text, exactly as authored. Force text number format first,
# then clear the style back to Normal so no extra formatting is left behind on the cell.
$priceUpdates = @{
    'D2' = '29.809.37'
    'D3' = '1.890.29'
    'D4' = '1.0000'
    'D5' = '0.7754'
    'D6' = '244.87'
    'D7' = '0.9999'
    'D8' = '0.3137'
    'D9' = '25.35'
    'D10' = '0.07231'
    'D11' = '0.08092'
    'D12' = '0.7665'
    'D13' = '1.940.94'
    'D14' = '5.486'
    'D15' = '92.37'
    'D16' = '6.193'
    'D17' = '29.821.61'
    'D18' = '13.92'
    'D19' = '243.14'
    'D20' = '0.000007773'
    'D21' = '2.161.80'
    'D22' = '1.001'
    'D23' = '8.181'
    'D24' = '1.000'
    'D25' = '0.1581'
    'D26' = '9.448'
    'D27' = '162.05'
    'D28' = '18.75'
    'D29' = '2.041'
    'D30' = '1.452'
    'D31' = '1.549'
    'D32' = '4.474'
    'D33' = '4.080'
    'D34' = '0.05519'
    'D35' = '1.256'
    'D36' = '0.7541'
    'D37' = '1.003'
    'D39' = '0.01922'
    'D40' = '2.777'
    'D41' = '1.156.60'
    'D42' = '73.80'
    'D43' = '0.4421'
    'D44' = '5.916'
    'D45' = '0.8474'
    'D46' = '0.9999'
    'D47' = '1.892'
    'D48' = '102.79'
    'D49' = '9.933'
    'D50' = '7.466'
    'D51' = '3.022'
}

foreach ($cell in $priceUpdates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $priceUpdates[$cell]
    $range.Style = "Normal"
}

Write-Output "Applied $($updates.Count + $priceUpdates.Count) cell updates"
